# "1st changes of mifos to finflux"
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# immediately before the old column N ("Late"), pushing the old
# N/O/P ("Late" / "heading" / "Outstanding") one column to the right
# (-> O/P/Q). The new column N picks up the column width that column M
# already had (Excel's normal "insert column" behavior of inheriting the
# format of the column to its left).
#
# Finally, the "Repayment schedule" tab is made the active sheet/tab
# (previously "Transactions" was active), with cell R7 selected there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Width (in the "characters" unit the ColumnWidth property uses) of the
# existing column M - the new column inherits this width, same as
# column M to its left.
$leftWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new blank column at N, shifting the old N:P right to O:Q.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet and select R7 on it (this
# also clears the previously active "Transactions" tab's tabSelected
# flag / updates the workbook's remembered active tab).
$ws.Activate() | Out-Null
$ws.Range("R7").Select() | Out-Null
